$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2310" and "_new" -> "_FV2404" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2310")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2404")
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used data range into a native Excel table (ListObject) ---
$rng = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
